# Applies the weekly update to the Pepino dulce (Feria Lagunitas de Puerto Montt) sheet.
# The underlying records (Calidad/Volumen/Precios) were reshuffled across a different
# set of dates (Fecha) for each row, so every affected cell is rewritten explicitly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 44251
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 20
$ws.Cells.Item(2, 11).Value = 18000
$ws.Cells.Item(2, 12).Value = 18000
$ws.Cells.Item(2, 13).Value = 18000
$ws.Cells.Item(2, 16).Value = 1000

$ws.Cells.Item(3, 4).Value = 44257
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 16000
$ws.Cells.Item(3, 12).Value = 16000
$ws.Cells.Item(3, 13).Value = 16000
$ws.Cells.Item(3, 16).Value = 889

$ws.Cells.Item(4, 4).Value = 44271
$ws.Cells.Item(4, 9).Value = "Especial"
$ws.Cells.Item(4, 10).Value = 70
$ws.Cells.Item(4, 11).Value = 18000
$ws.Cells.Item(4, 12).Value = 18000
$ws.Cells.Item(4, 13).Value = 18000
$ws.Cells.Item(4, 16).Value = 1000

$ws.Cells.Item(5, 4).Value = 44260
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 70
$ws.Cells.Item(5, 11).Value = 16000
$ws.Cells.Item(5, 12).Value = 16000
$ws.Cells.Item(5, 13).Value = 16000
$ws.Cells.Item(5, 16).Value = 889

$ws.Cells.Item(6, 4).Value = 44264
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 80
$ws.Cells.Item(6, 11).Value = 16000
$ws.Cells.Item(6, 12).Value = 16000
$ws.Cells.Item(6, 13).Value = 16000
$ws.Cells.Item(6, 16).Value = 889

$ws.Cells.Item(7, 4).Value = 44301
$ws.Cells.Item(7, 9).Value = "Especial"
$ws.Cells.Item(7, 10).Value = 30
$ws.Cells.Item(7, 11).Value = 18000
$ws.Cells.Item(7, 12).Value = 18000
$ws.Cells.Item(7, 13).Value = 18000
$ws.Cells.Item(7, 16).Value = 1000

$ws.Cells.Item(8, 4).Value = 44253
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 70
$ws.Cells.Item(8, 11).Value = 18000
$ws.Cells.Item(8, 12).Value = 18000
$ws.Cells.Item(8, 13).Value = 18000
$ws.Cells.Item(8, 16).Value = 1000

$ws.Cells.Item(9, 4).Value = 44238
$ws.Cells.Item(9, 9).Value = "Especial"
$ws.Cells.Item(9, 10).Value = 50
$ws.Cells.Item(9, 11).Value = 20000
$ws.Cells.Item(9, 12).Value = 20000
$ws.Cells.Item(9, 13).Value = 20000
$ws.Cells.Item(9, 16).Value = 1111

$ws.Cells.Item(10, 4).Value = 44364
$ws.Cells.Item(10, 9).Value = "Especial"
$ws.Cells.Item(10, 10).Value = 30
$ws.Cells.Item(10, 11).Value = 20000
$ws.Cells.Item(10, 12).Value = 20000
$ws.Cells.Item(10, 13).Value = 20000
$ws.Cells.Item(10, 16).Value = 1111

$ws.Cells.Item(11, 4).Value = 44302
$ws.Cells.Item(11, 9).Value = "Especial"
$ws.Cells.Item(11, 10).Value = 70
$ws.Cells.Item(11, 11).Value = 18000
$ws.Cells.Item(11, 12).Value = 18000
$ws.Cells.Item(11, 13).Value = 18000
$ws.Cells.Item(11, 16).Value = 1000

$ws.Cells.Item(12, 4).Value = 44320
$ws.Cells.Item(12, 9).Value = "Especial"
$ws.Cells.Item(12, 10).Value = 90
$ws.Cells.Item(12, 11).Value = 17500
$ws.Cells.Item(12, 12).Value = 17500
$ws.Cells.Item(12, 13).Value = 17500
$ws.Cells.Item(12, 16).Value = 972

$ws.Cells.Item(13, 4).Value = 44365
$ws.Cells.Item(13, 9).Value = "Especial"
$ws.Cells.Item(13, 10).Value = 50
$ws.Cells.Item(13, 11).Value = 20000
$ws.Cells.Item(13, 12).Value = 20000
$ws.Cells.Item(13, 13).Value = 20000
$ws.Cells.Item(13, 16).Value = 1111

$ws.Cells.Item(14, 4).Value = 44312
$ws.Cells.Item(14, 9).Value = "Especial"
$ws.Cells.Item(14, 10).Value = 15
$ws.Cells.Item(14, 11).Value = 18000
$ws.Cells.Item(14, 12).Value = 18000
$ws.Cells.Item(14, 13).Value = 18000
$ws.Cells.Item(14, 16).Value = 1000

$ws.Cells.Item(15, 4).Value = 44309
$ws.Cells.Item(15, 9).Value = "Especial"
$ws.Cells.Item(15, 10).Value = 80
$ws.Cells.Item(15, 11).Value = 18000
$ws.Cells.Item(15, 12).Value = 18000
$ws.Cells.Item(15, 13).Value = 18000
$ws.Cells.Item(15, 16).Value = 1000

$ws.Cells.Item(16, 4).Value = 44323
$ws.Cells.Item(16, 9).Value = "Especial"
$ws.Cells.Item(16, 10).Value = 70
$ws.Cells.Item(16, 11).Value = 18000
$ws.Cells.Item(16, 12).Value = 18000
$ws.Cells.Item(16, 13).Value = 18000
$ws.Cells.Item(16, 16).Value = 1000

$ws.Cells.Item(17, 4).Value = 44295
$ws.Cells.Item(17, 9).Value = "Especial"
$ws.Cells.Item(17, 10).Value = 80
$ws.Cells.Item(17, 11).Value = 16000
$ws.Cells.Item(17, 12).Value = 16000
$ws.Cells.Item(17, 13).Value = 16000
$ws.Cells.Item(17, 16).Value = 889

$ws.Cells.Item(18, 4).Value = 44246
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 60
$ws.Cells.Item(18, 11).Value = 18000
$ws.Cells.Item(18, 12).Value = 18000
$ws.Cells.Item(18, 13).Value = 18000
$ws.Cells.Item(18, 16).Value = 1000

$ws.Cells.Item(19, 4).Value = 44327
$ws.Cells.Item(19, 9).Value = "Especial"
$ws.Cells.Item(19, 10).Value = 80
$ws.Cells.Item(19, 11).Value = 16000
$ws.Cells.Item(19, 12).Value = 16000
$ws.Cells.Item(19, 13).Value = 16000
$ws.Cells.Item(19, 16).Value = 889

$ws.Cells.Item(20, 4).Value = 44232
$ws.Cells.Item(20, 9).Value = "Especial"
$ws.Cells.Item(20, 10).Value = 50
$ws.Cells.Item(20, 11).Value = 22000
$ws.Cells.Item(20, 12).Value = 22000
$ws.Cells.Item(20, 13).Value = 22000
$ws.Cells.Item(20, 16).Value = 1222

$ws.Cells.Item(21, 4).Value = 44326
$ws.Cells.Item(21, 9).Value = "Especial"
$ws.Cells.Item(21, 10).Value = 15
$ws.Cells.Item(21, 11).Value = 18000
$ws.Cells.Item(21, 12).Value = 18000
$ws.Cells.Item(21, 13).Value = 18000
$ws.Cells.Item(21, 16).Value = 1000

$ws.Cells.Item(22, 4).Value = 44292
$ws.Cells.Item(22, 9).Value = "Especial"
$ws.Cells.Item(22, 10).Value = 70
$ws.Cells.Item(22, 11).Value = 17000
$ws.Cells.Item(22, 12).Value = 17000
$ws.Cells.Item(22, 13).Value = 17000
$ws.Cells.Item(22, 16).Value = 944

$ws.Cells.Item(23, 4).Value = 44252
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 40
$ws.Cells.Item(23, 11).Value = 18000
$ws.Cells.Item(23, 12).Value = 18000
$ws.Cells.Item(23, 13).Value = 18000
$ws.Cells.Item(23, 16).Value = 1000

$ws.Cells.Item(24, 4).Value = 44245
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 40
$ws.Cells.Item(24, 11).Value = 18000
$ws.Cells.Item(24, 12).Value = 18000
$ws.Cells.Item(24, 13).Value = 18000
$ws.Cells.Item(24, 16).Value = 1000

$ws.Cells.Item(25, 4).Value = 44243
$ws.Cells.Item(25, 9).Value = "Especial"
$ws.Cells.Item(25, 10).Value = 60
$ws.Cells.Item(25, 11).Value = 20000
$ws.Cells.Item(25, 12).Value = 20000
$ws.Cells.Item(25, 13).Value = 20000
$ws.Cells.Item(25, 16).Value = 1111

$ws.Cells.Item(26, 4).Value = 44250
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 60
$ws.Cells.Item(26, 11).Value = 18000
$ws.Cells.Item(26, 12).Value = 18000
$ws.Cells.Item(26, 13).Value = 18000
$ws.Cells.Item(26, 16).Value = 1000

$ws.Cells.Item(27, 4).Value = 44236
$ws.Cells.Item(27, 9).Value = "Especial"
$ws.Cells.Item(27, 10).Value = 60
$ws.Cells.Item(27, 11).Value = 20000
$ws.Cells.Item(27, 12).Value = 20000
$ws.Cells.Item(27, 13).Value = 20000
$ws.Cells.Item(27, 16).Value = 1111

$ws.Cells.Item(28, 4).Value = 44306
$ws.Cells.Item(28, 9).Value = "Especial"
$ws.Cells.Item(28, 10).Value = 80
$ws.Cells.Item(28, 11).Value = 18000
$ws.Cells.Item(28, 12).Value = 18000
$ws.Cells.Item(28, 13).Value = 18000
$ws.Cells.Item(28, 16).Value = 1000

$ws.Cells.Item(29, 4).Value = 44285
$ws.Cells.Item(29, 9).Value = "Especial"
$ws.Cells.Item(29, 10).Value = 70
$ws.Cells.Item(29, 11).Value = 18000
$ws.Cells.Item(29, 12).Value = 18000
$ws.Cells.Item(29, 13).Value = 18000
$ws.Cells.Item(29, 16).Value = 1000

$ws.Cells.Item(30, 4).Value = 44371
$ws.Cells.Item(30, 9).Value = "Especial"
$ws.Cells.Item(30, 10).Value = 20
$ws.Cells.Item(30, 11).Value = 20000
$ws.Cells.Item(30, 12).Value = 20000
$ws.Cells.Item(30, 13).Value = 20000
$ws.Cells.Item(30, 16).Value = 1111

$ws.Cells.Item(31, 4).Value = 44242
$ws.Cells.Item(31, 9).Value = "Especial"
$ws.Cells.Item(31, 10).Value = 50
$ws.Cells.Item(31, 11).Value = 20000
$ws.Cells.Item(31, 12).Value = 20000
$ws.Cells.Item(31, 13).Value = 20000
$ws.Cells.Item(31, 16).Value = 1111

$ws.Cells.Item(32, 4).Value = 44313
$ws.Cells.Item(32, 9).Value = "Especial"
$ws.Cells.Item(32, 10).Value = 80
$ws.Cells.Item(32, 11).Value = 18000
$ws.Cells.Item(32, 12).Value = 18000
$ws.Cells.Item(32, 13).Value = 18000
$ws.Cells.Item(32, 16).Value = 1000

$ws.Cells.Item(33, 4).Value = 44357
$ws.Cells.Item(33, 9).Value = "Especial"
$ws.Cells.Item(33, 10).Value = 15
$ws.Cells.Item(33, 11).Value = 20000
$ws.Cells.Item(33, 12).Value = 20000
$ws.Cells.Item(33, 13).Value = 20000
$ws.Cells.Item(33, 16).Value = 1111

$ws.Cells.Item(34, 4).Value = 44259
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 70
$ws.Cells.Item(34, 11).Value = 16000
$ws.Cells.Item(34, 12).Value = 16000
$ws.Cells.Item(34, 13).Value = 16000
$ws.Cells.Item(34, 16).Value = 889

$ws.Cells.Item(35, 4).Value = 44278
$ws.Cells.Item(35, 9).Value = "Especial"
$ws.Cells.Item(35, 10).Value = 70
$ws.Cells.Item(35, 11).Value = 18000
$ws.Cells.Item(35, 12).Value = 18000
$ws.Cells.Item(35, 13).Value = 18000
$ws.Cells.Item(35, 16).Value = 1000

$ws.Cells.Item(36, 4).Value = 44280
$ws.Cells.Item(36, 9).Value = "Especial"
$ws.Cells.Item(36, 10).Value = 40
$ws.Cells.Item(36, 11).Value = 18000
$ws.Cells.Item(36, 12).Value = 18000
$ws.Cells.Item(36, 13).Value = 18000
$ws.Cells.Item(36, 16).Value = 1000

$ws.Cells.Item(37, 4).Value = 44267
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 70
$ws.Cells.Item(37, 11).Value = 16000
$ws.Cells.Item(37, 12).Value = 16000
$ws.Cells.Item(37, 13).Value = 16000
$ws.Cells.Item(37, 16).Value = 889

$ws.Cells.Item(38, 4).Value = 44270
$ws.Cells.Item(38, 9).Value = "Primera"
$ws.Cells.Item(38, 10).Value = 15
$ws.Cells.Item(38, 11).Value = 16000
$ws.Cells.Item(38, 12).Value = 16000
$ws.Cells.Item(38, 13).Value = 16000
$ws.Cells.Item(38, 16).Value = 889

$ws.Cells.Item(39, 4).Value = 44299
$ws.Cells.Item(39, 9).Value = "Especial"
$ws.Cells.Item(39, 10).Value = 80
$ws.Cells.Item(39, 11).Value = 18000
$ws.Cells.Item(39, 12).Value = 18000
$ws.Cells.Item(39, 13).Value = 18000
$ws.Cells.Item(39, 16).Value = 1000

$ws.Cells.Item(40, 4).Value = 44316
$ws.Cells.Item(40, 9).Value = "Especial"
$ws.Cells.Item(40, 10).Value = 70
$ws.Cells.Item(40, 11).Value = 18000
$ws.Cells.Item(40, 12).Value = 18000
$ws.Cells.Item(40, 13).Value = 18000
$ws.Cells.Item(40, 16).Value = 1000

$ws.Cells.Item(41, 4).Value = 44239
$ws.Cells.Item(41, 9).Value = "Especial"
$ws.Cells.Item(41, 10).Value = 60
$ws.Cells.Item(41, 11).Value = 20000
$ws.Cells.Item(41, 12).Value = 20000
$ws.Cells.Item(41, 13).Value = 20000
$ws.Cells.Item(41, 16).Value = 1111

$ws.Cells.Item(42, 4).Value = 44350
$ws.Cells.Item(42, 9).Value = "Primera"
$ws.Cells.Item(42, 10).Value = 20
$ws.Cells.Item(42, 11).Value = 20000
$ws.Cells.Item(42, 12).Value = 20000
$ws.Cells.Item(42, 13).Value = 20000
$ws.Cells.Item(42, 16).Value = 1111

$ws.Cells.Item(43, 4).Value = 44274
$ws.Cells.Item(43, 9).Value = "Primera"
$ws.Cells.Item(43, 10).Value = 70
$ws.Cells.Item(43, 11).Value = 16000
$ws.Cells.Item(43, 12).Value = 16000
$ws.Cells.Item(43, 13).Value = 16000
$ws.Cells.Item(43, 16).Value = 889

